$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) is always treated as text, matching the
# source data which stores prices as formatted strings (e.g. "42.184.89").
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '42.184.89'
$ws.Range("E2").Value = '  -2.18%  '

# Row 3
$ws.Range("D3").Value = '2.298.49'
$ws.Range("E3").Value = '  -3.26%  '

# Row 4
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").Value = '317.81'
$ws.Range("E5").Value = '  +0.33%  '

# Row 6
$ws.Range("D6").Value = '104.44'
$ws.Range("E6").Value = '  -4.48%  '

# Row 7
$ws.Range("D7").Value = '0.626'
$ws.Range("E7").Value = '  -2.14%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").Value = '0.610'
$ws.Range("E9").Value = '  -2.15%  '

# Row 10
$ws.Range("D10").Value = '39.62'
$ws.Range("E10").Value = '  -3.94%  '

# Row 11
$ws.Range("E11").Value = '  -2.52%  '

# Row 12
$ws.Range("D12").Value = '8.33'
$ws.Range("E12").Value = '  -3.64%  '

# Row 13
$ws.Range("E13").Value = '  -0.47%  '

# Row 14
$ws.Range("D14").Value = '0.965'
$ws.Range("E14").Value = '  -4.78%  '

# Row 15
$ws.Range("D15").Value = '15.33'
$ws.Range("E15").Value = '  -4.53%  '

# Row 16
$ws.Range("D16").Value = '2.645.73'
$ws.Range("E16").Value = '  -3.31%  '

# Row 17
$ws.Range("D17").Value = '2.286.80'
$ws.Range("E17").Value = '  -4.88%  '

# Row 18
$ws.Range("D18").Value = '42.261.26'
$ws.Range("E18").Value = '  -2.00%  '

# Row 19
$ws.Range("D19").Value = '7.38'
$ws.Range("E19").Value = '  -4.55%  '

# Row 20
$ws.Range("D20").Value = '0.0000106'
$ws.Range("E20").Value = '  -1.27%  '

# Row 21
$ws.Range("D21").Value = '73.46'
$ws.Range("E21").Value = '  -4.07%  '

# Row 22
$ws.Range("D22").Value = '3.62'
$ws.Range("E22").Value = '  +0.21%  '

# Row 23
$ws.Range("D23").Value = '278.96'
$ws.Range("E23").Value = '  +3.47%  '

# Row 24
$ws.Range("D24").Value = '10.76'
$ws.Range("E24").Value = '  +11.95%  '

# Row 25
$ws.Range("E25").Value = '  -3.33%  '

# Row 26
$ws.Range("E26").Value = '  +0.20%  '

# Row 27
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = '2.43'
$ws.Range("E27").Value = '  +7.64%  '

# Row 28
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").Value = '10.86'
$ws.Range("E28").Value = '  -5.74%  '

# Row 29
$ws.Range("E29").Value = '  -3.14%  '

# Row 30
$ws.Range("D30").Value = '36.09'
$ws.Range("E30").Value = '  -3.41%  '

# Row 31
$ws.Range("D31").Value = '163.68'
$ws.Range("E31").Value = '  -3.19%  '

# Row 32
$ws.Range("D32").Value = '0.0873'
$ws.Range("E32").Value = '  -4.20%  '

# Row 33
$ws.Range("D33").Value = '5.84'
$ws.Range("E33").Value = '  -5.74%  '

# Row 34
$ws.Range("E34").Value = '  -5.68%  '

# Row 35
$ws.Range("D35").Value = '0.137'
$ws.Range("E35").Value = '  +3.78%  '

# Row 36
$ws.Range("E36").Value = '  -6.47%  '

# Row 37
$ws.Range("D37").Value = '4.60'
$ws.Range("E37").Value = '  -3.02%  '

# Row 38
$ws.Range("D38").Value = '0.0349'
$ws.Range("E38").Value = '  -4.05%  '

# Row 39
$ws.Range("D39").Value = '3.75'
$ws.Range("E39").Value = '  -3.32%  '

# Row 40
$ws.Range("E40").Value = '  +2.69%  '

# Row 41
$ws.Range("D41").Value = '100.04'
$ws.Range("E41").Value = '  -4.88%  '

# Row 42
$ws.Range("E42").Value = '  -5.01%  '

# Row 43
$ws.Range("D43").Value = '69.48'
$ws.Range("E43").Value = '  -3.48%  '

# Row 44
$ws.Range("D44").Value = '0.226'
$ws.Range("E44").Value = '  -5.41%  '

# Row 45
$ws.Range("E45").Value = '  +0.08%  '

# Row 46
$ws.Range("D46").Value = '12.06'
$ws.Range("E46").Value = '  -5.86%  '

# Row 47
$ws.Range("D47").Value = '112.12'
$ws.Range("E47").Value = '  -2.47%  '

# Row 48
$ws.Range("D48").Value = '77.33'
$ws.Range("E48").Value = '  -3.93%  '

# Row 49
$ws.Range("D49").Value = '8.93'
$ws.Range("E49").Value = '  -3.17%  '

# Row 50
$ws.Range("D50").Value = '5.29'
$ws.Range("E50").Value = '  -5.48%  '

# Row 51
$ws.Range("D51").Value = '1.606.25'
$ws.Range("E51").Value = '  +1.32%  '
